$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The player roster table (rows 2-19) has been re-sorted into a new row
# order. Header row 1 (Oyuncu Adı / Pozisyon / Takım) is unchanged.
$data = @(
    @("Chris Paul",          "PG",      "San Antonio Spurs"),
    @("Russell Westbrook",   "PG,SG",   "Denver Nuggets"),
    @("Stephon Castle",      "PG,SG",   "San Antonio Spurs"),
    @("Cole Anthony",        "PG",      "Orlando Magic"),
    @("Jaylen Brown",        "SG,SF",   "Boston Celtics"),
    @("Paolo Banchero",      "SF,PF",   "Orlando Magic"),
    @("Pascal Siakam",       "SF,PF,C", "Indiana Pacers"),
    @("Bol Bol",             "PF,C",    "Phoenix Suns"),
    @("Nikola Jokic",        "C",       "Denver Nuggets"),
    @("Jakob Poeltl",        "C",       "Toronto Raptors"),
    @("Payton Pritchard",    "PG,SG",   "Boston Celtics"),
    @("Jerami Grant",        "SF,PF",   "Portland Trail Blazers"),
    @("Chet Holmgren",       "PF,C",    "Oklahoma City Thunder"),
    @("Jalen Green",         "PG,SG",   "Houston Rockets"),
    @("Brandin Podziemski",  "PG,SG",   "Golden State Warriors"),
    @("Rudy Gobert",         "C",       "Minnesota Timberwolves"),
    @("Jalen Suggs",         "PG,SG",   "Orlando Magic"),
    @("Deni Avdija",         "SF,PF",   "Portland Trail Blazers")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}
